$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 138, shifting existing rows 138:242 down to 139:243
$ws.Rows(138).Insert()

# Populate the newly inserted row 138 with the new data record
$ws.Cells.Item(138,1).Value = 9
$ws.Cells.Item(138,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(138,3).Value = "Metropolitana"
$ws.Cells.Item(138,4).Value = 44762
$ws.Cells.Item(138,5).Value = 13
$ws.Cells.Item(138,6).Value = 100112026
$ws.Cells.Item(138,7).Value = "Haba"
$ws.Cells.Item(138,8).Value = "Sin especificar"
$ws.Cells.Item(138,9).Value = "Primera"
$ws.Cells.Item(138,10).Value = 52
$ws.Cells.Item(138,11).Value = 16000
$ws.Cells.Item(138,12).Value = 16000
$ws.Cells.Item(138,13).Value = 16000
$ws.Cells.Item(138,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(138,15).Value = "Región de Coquimbo"
$ws.Cells.Item(138,16).Value = 640
$ws.Cells.Item(138,17).Value = 25
$ws.Cells.Item(138,18).Value = "Hortaliza"
